$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the previously-missing calculation for row 5 (PriceChange / UpDown) ---
$ws.Cells.Item(5, 24).Value = 0.11999500000000296   # X5 - PriceChange
$ws.Cells.Item(5, 25).Value = "Up"                  # Y5 - UpDown

# --- Append a new scan result as row 6 ---
$ws.Cells.Item(6, 1).Value = 42647.887013888889     # A6 - Date
$ws.Cells.Item(6, 1).NumberFormat = "m/d/yy h:mm"

$ws.Cells.Item(6, 2).Value = -20                    # B6 - ScoreFinal
$ws.Cells.Item(6, 3).Value = "Strong Sell"          # C6 - Verdict
$ws.Cells.Item(6, 4).Value = 0                      # D6 - totalSentiment
$ws.Cells.Item(6, 5).Value = 0                      # E6 - wordCount
$ws.Cells.Item(6, 6).Value = 0                      # F6 - sentenceCount
$ws.Cells.Item(6, 7).Value = 0                      # G6 - posWordPercentage
$ws.Cells.Item(6, 8).Value = 0                      # H6 - negWordPercentage
$ws.Cells.Item(6, 9).Value = 0                      # I6 - posPhrasePercentage
$ws.Cells.Item(6, 10).Value = 0                     # J6 - negPhrasePercentage
$ws.Cells.Item(6, 11).Value = 0                     # K6 - ElapsedMs
$ws.Cells.Item(6, 12).Value = 0                     # L6 - posWordCount
$ws.Cells.Item(6, 13).Value = 0                     # M6 - negWordCount
$ws.Cells.Item(6, 14).Value = 0                     # N6 - positivePhraseCount
$ws.Cells.Item(6, 15).Value = 0                     # O6 - negativePhraseCount
$ws.Cells.Item(6, 16).Value = "Random"              # P6 - Method
$ws.Cells.Item(6, 17).Value = 42.459412013272512    # Q6 - RSI
$ws.Cells.Item(6, 18).Value = 0                     # R6 - PEG

$ws.Cells.Item(6, 19).Value = -0.0112               # S6 - 200Moving%
$ws.Cells.Item(6, 19).NumberFormat = "0.00%"
$ws.Cells.Item(6, 20).Value = -0.0367               # T6 - 50Moving%
$ws.Cells.Item(6, 20).NumberFormat = "0.00%"

$ws.Cells.Item(6, 21).Value = 14.56                 # U6 - PriceBook
$ws.Cells.Item(6, 22).Value = "N/A"                 # V6 - Dividend
$ws.Cells.Item(6, 23).Value = -2                    # W6 - Bollinger
